$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1, Col 1: "50÷9=5, 5" -> "56÷5=11, 1"
$cell = $t.Cell(1, 1)
$cell.Range.Text = "56÷5=11, 1"

# Row 1, Col 2: "49÷6=8, 1" -> "85÷9=9, 4"
$cell = $t.Cell(1, 2)
$cell.Range.Text = "85÷9=9, 4"

# Row 1, Col 3: "23÷4=5, 3" -> "41÷6=6, 5"
$cell = $t.Cell(1, 3)
$cell.Range.Text = "41÷6=6, 5"

# Row 1, Col 4: "63÷2=31, 1" -> "14÷2=7, 0"
$cell = $t.Cell(1, 4)
$cell.Range.Text = "14÷2=7, 0"

# Row 1, Col 5: "48÷6=8, 0" -> "23÷2=11, 1"
$cell = $t.Cell(1, 5)
$cell.Range.Text = "23÷2=11, 1"

# Row 5, Col 1: "87÷3=29, 0" -> "36÷8=4, 4"
$cell = $t.Cell(5, 1)
$cell.Range.Text = "36÷8=4, 4"

# Row 5, Col 2: "95÷3=31, 2" -> "63÷8=7, 7"
$cell = $t.Cell(5, 2)
$cell.Range.Text = "63÷8=7, 7"

# Row 5, Col 3: "64÷2=32, 0" -> "59÷7=8, 3"
$cell = $t.Cell(5, 3)
$cell.Range.Text = "59÷7=8, 3"

# Row 5, Col 4: "93÷8=11, 5" -> "63÷4=15, 3"
$cell = $t.Cell(5, 4)
$cell.Range.Text = "63÷4=15, 3"

# Row 5, Col 5: "32÷9=3, 5" -> "59÷4=14, 3"
$cell = $t.Cell(5, 5)
$cell.Range.Text = "59÷4=14, 3"

# Row 9, Col 1: "81÷3=27, 0" -> "22÷9=2, 4"
$cell = $t.Cell(9, 1)
$cell.Range.Text = "22÷9=2, 4"

# Row 9, Col 2: "23÷9=2, 5" -> "65÷8=8, 1"
$cell = $t.Cell(9, 2)
$cell.Range.Text = "65÷8=8, 1"

# Row 9, Col 3: "64÷2=32, 0" -> "89÷5=17, 4"
$cell = $t.Cell(9, 3)
$cell.Range.Text = "89÷5=17, 4"

# Row 9, Col 4: "50÷7=7, 1" -> "78÷3=26, 0"
$cell = $t.Cell(9, 4)
$cell.Range.Text = "78÷3=26, 0"

# Row 9, Col 5: "74÷5=14, 4" -> "86÷6=14, 2"
$cell = $t.Cell(9, 5)
$cell.Range.Text = "86÷6=14, 2"

# Row 13, Col 1: "86÷8=10, 6" -> "37÷3=12, 1"
$cell = $t.Cell(13, 1)
$cell.Range.Text = "37÷3=12, 1"

# Row 13, Col 2: "65÷5=13, 0" -> "87÷4=21, 3"
$cell = $t.Cell(13, 2)
$cell.Range.Text = "87÷4=21, 3"

# Row 13, Col 3: "14÷2=7, 0" -> "25÷3=8, 1"
$cell = $t.Cell(13, 3)
$cell.Range.Text = "25÷3=8, 1"

# Row 13, Col 4: "13÷2=6, 1" -> "50÷4=12, 2"
$cell = $t.Cell(13, 4)
$cell.Range.Text = "50÷4=12, 2"

# Row 13, Col 5: "74÷3=24, 2" -> "91÷8=11, 3"
$cell = $t.Cell(13, 5)
$cell.Range.Text = "91÷8=11, 3"

# Row 17, Col 1: "85÷9=9, 4" -> "30÷3=10, 0"
$cell = $t.Cell(17, 1)
$cell.Range.Text = "30÷3=10, 0"

# Row 17, Col 2: "34÷8=4, 2" -> "20÷2=10, 0"
$cell = $t.Cell(17, 2)
$cell.Range.Text = "20÷2=10, 0"

# Row 17, Col 3: "91÷2=45, 1" -> "95÷9=10, 5"
$cell = $t.Cell(17, 3)
$cell.Range.Text = "95÷9=10, 5"

# Row 17, Col 4: "40÷9=4, 4" -> "39÷5=7, 4"
$cell = $t.Cell(17, 4)
$cell.Range.Text = "39÷5=7, 4"

# Row 17, Col 5: "75÷9=8, 3" -> "76÷7=10, 6"
$cell = $t.Cell(17, 5)
$cell.Range.Text = "76÷7=10, 6"
